$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Before: A1=State B1=Subject C1=Grade D1=Year(s)
# After:  A1=FIPS  B1=State   C1=Subject D1=Grade(s) E1=Year(s) F1=Notes
$ws.Range("A1").Value = "FIPS"
$ws.Range("B1").Value = "State"
$ws.Range("C1").Value = "Subject"
$ws.Range("D1").Value = "Grade(s)"
$ws.Range("E1").Value = "Year(s)"
$ws.Range("F1").Value = "Notes"

# --- Data rows (2-9) ---
$data = @(
    @(5,  "Arkansas",  "Math", 8,            "2009-2010", $null),
    @(22, "Louisiana", "Math", "3, 4",        2018,        $null),
    @(22, "Louisiana", "RLA",  "3, 4",        2018,        $null),
    @(29, "Missouri",  "Math", 8,            "2013-2018", $null),
    @(36, "New York",  "Math", "6, 7, 8",     2014,        "Revisit x-axis; ends in 2014"),
    @(36, "New York",  "RLA",  "6, 7, 8",     2014,        "Revisit x-axis; ends in 2014"),
    @(48, "Texas",     "Math", "7, 8",        "2012-2018", $null),
    @(51, "Virginia",  "Math", "5, 6, 7, 8",  "2009-2018", "Revisit x-axis")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    if ($entry[5] -ne $null) {
        $ws.Cells.Item($row, 6).Value = $entry[5]
    }
    $row = $row + 1
}

# --- Column widths ---
# Target OOXML width is 25.83203125; this runtime quantizes ColumnWidth to
# 1/6-character steps, so 25 is the input that round-trips closest to it.
$ws.Columns.Item(6).ColumnWidth = 25

# --- Selection ---
$ws.Range("F13").Select()
